# Update generated output numbers (view/participation counts) across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 240
$ws1.Cells.Item(5, 6).Value  = 9368
$ws1.Cells.Item(6, 6).Value  = 9368
$ws1.Cells.Item(7, 6).Value  = 584
$ws1.Cells.Item(10, 6).Value = 266
$ws1.Cells.Item(11, 6).Value = 403
$ws1.Cells.Item(13, 6).Value = 168
$ws1.Cells.Item(15, 6).Value = 11983
$ws1.Cells.Item(16, 6).Value = 11983
$ws1.Cells.Item(27, 6).Value = 173
$ws1.Cells.Item(29, 6).Value = 2721
$ws1.Cells.Item(31, 6).Value = 100
$ws1.Cells.Item(32, 6).Value = 2097
$ws1.Cells.Item(37, 6).Value = 993
$ws1.Cells.Item(38, 6).Value = 4190
$ws1.Cells.Item(39, 6).Value = 3621
$ws1.Cells.Item(40, 6).Value = 505
$ws1.Cells.Item(41, 6).Value = 2620
$ws1.Cells.Item(43, 6).Value = 1314
$ws1.Cells.Item(44, 6).Value = 192
$ws1.Cells.Item(46, 6).Value = 411
$ws1.Cells.Item(47, 6).Value = 504
$ws1.Cells.Item(48, 6).Value = 64
$ws1.Cells.Item(49, 6).Value = 213
$ws1.Cells.Item(50, 6).Value = 124
$ws1.Cells.Item(51, 6).Value = 133

# Sheet "演出" (Performances) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(16, 6).Value = 6
$ws2.Cells.Item(18, 6).Value = 8
$ws2.Cells.Item(19, 6).Value = 185
$ws2.Cells.Item(21, 6).Value = 35

# Sheet "全部类型" (All types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(8, 6).Value  = 240
$ws4.Cells.Item(9, 6).Value  = 9368
$ws4.Cells.Item(10, 6).Value = 584
$ws4.Cells.Item(13, 6).Value = 266
$ws4.Cells.Item(14, 6).Value = 403
$ws4.Cells.Item(16, 6).Value = 168
$ws4.Cells.Item(17, 6).Value = 11983
$ws4.Cells.Item(18, 6).Value = 11983
$ws4.Cells.Item(28, 6).Value = 173
$ws4.Cells.Item(30, 6).Value = 2721
$ws4.Cells.Item(32, 6).Value = 100
$ws4.Cells.Item(33, 6).Value = 2097
$ws4.Cells.Item(37, 6).Value = 8
$ws4.Cells.Item(39, 6).Value = 993
$ws4.Cells.Item(40, 6).Value = 185
$ws4.Cells.Item(42, 6).Value = 3621
$ws4.Cells.Item(45, 6).Value = 1314
$ws4.Cells.Item(46, 6).Value = 192
$ws4.Cells.Item(47, 6).Value = 411
$ws4.Cells.Item(49, 6).Value = 504
$ws4.Cells.Item(50, 6).Value = 64
$ws4.Cells.Item(51, 6).Value = 213
